$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Enterprises density (per 1000 people)" row (row 13): Micro / SMEs / MSMEs
$ws.Range("B13").Value = "'40.25"
$ws.Range("C13").Value = "'0.69"
$ws.Range("D13").Value = "'40.94"

# "Employment (% of total)" row (row 14): Micro / SMEs / MSMEs
$ws.Range("B14").Value = "'76.94"
$ws.Range("C14").Value = "'12.67"
$ws.Range("D14").Value = "'89.61"

# "Enterprises (% of total)" row (row 16): SMEs / MSMEs (Micro, B16, stays 98.1)
$ws.Range("C16").Value = "'1.67"
$ws.Range("D16").Value = "'99.78"
